$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 224.8
$ws.Range("I2").Value = 71.77778000000001
$ws.Range("K2").Value = 71.77778000000001
$ws.Range("M2").Value = 41.22221999999999
$ws.Range("H13").Value = 3500
$ws.Range("J13").Value = 3500
$ws.Range("L13").Value = 3500
$ws.Range("N13").Value = -3838
$ws.Range("H38").Value = 277.5
$ws.Range("I38").Value = 277.5
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 832.5
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -460.5
$ws.Range("N38").ClearContents()
$ws.Range("H69").Value = 16460.072
$ws.Range("I69").Value = 8200
$ws.Range("J69").Value = 21049
$ws.Range("K69").Value = 24600
$ws.Range("L69").Value = 63147
$ws.Range("M69").Value = -23726
$ws.Range("N69").Value = -64895
$ws.Range("H72").Value = 16460.072
$ws.Range("I72").Value = 8200
$ws.Range("J72").Value = 21049
$ws.Range("K72").Value = 73800
$ws.Range("L72").Value = 189441
$ws.Range("M72").Value = -69432
$ws.Range("N72").Value = -198177
$ws.Range("H137").Value = 6785.4443
$ws.Range("I137").Value = 6739.231
$ws.Range("K137").Value = 20217.693
$ws.Range("M137").Value = -17667.693

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 38750
$ws.Range("J24").Value = 38750
$ws.Range("L24").Value = 38750
$ws.Range("N24").Value = -39498
$ws.Range("H32").Value = 1773.1143
$ws.Range("I32").Value = 1773.1143
$ws.Range("K32").Value = 1773.1143
$ws.Range("M32").Value = -1486.1143
$ws.Range("H97").Value = 2399.8333
$ws.Range("I97").Value = 879.8
$ws.Range("K97").Value = 879.8
$ws.Range("M97").Value = -383.8
$ws.Range("H100").Value = 38750
$ws.Range("J100").Value = 38750
$ws.Range("L100").Value = 38750
$ws.Range("N100").Value = -40914
$ws.Range("H101").Value = 35159.8
$ws.Range("J101").Value = 35159.8
$ws.Range("L101").Value = 35159.8
$ws.Range("N101").Value = -41649.8
$ws.Range("H132").Value = 3312.3428
$ws.Range("I132").Value = 3387.9033
$ws.Range("K132").Value = 10163.7099
$ws.Range("M132").Value = -7633.7099

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1701.6111
$ws.Range("I20").Value = 1604.7273
$ws.Range("K20").Value = 1604.7273
$ws.Range("M20").Value = -1357.7273
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H36").Value = 1250
$ws.Range("I36").Value = 1000
$ws.Range("K36").Value = 1000
$ws.Range("M36").Value = -466
$ws.Range("H80").Value = 6695.091
$ws.Range("I80").Value = 110.75
$ws.Range("K80").Value = 110.75
$ws.Range("M80").Value = 887.25
$ws.Range("H83").Value = 6695.091
$ws.Range("I83").Value = 110.75
$ws.Range("K83").Value = 553.75
$ws.Range("M83").Value = 4438.25
$ws.Range("H86").Value = 7298.2
$ws.Range("I86").Value = 3426.2856
$ws.Range("K86").Value = 3426.2856
$ws.Range("M86").Value = -2303.2856
$ws.Range("H89").Value = 7298.2
$ws.Range("I89").Value = 3426.2856
$ws.Range("K89").Value = 17131.428
$ws.Range("M89").Value = -11515.428
$ws.Range("H100").Value = 23999.8
$ws.Range("J100").Value = 23999.8
$ws.Range("L100").Value = 23999.8
$ws.Range("N100").Value = -26163.8
$ws.Range("H105").Value = 6262.3125
$ws.Range("I105").Value = 5976.6924
$ws.Range("K105").Value = 5976.6924
$ws.Range("M105").Value = -4229.6924

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 25000
$ws.Range("J25").Value = 25000
$ws.Range("L25").Value = 25000
$ws.Range("N25").Value = -25348
$ws.Range("H37").Value = 25000
$ws.Range("I37").Value = 25000
$ws.Range("K37").Value = 25000
$ws.Range("M37").Value = -24893
$ws.Range("H60").Value = 53401.668
$ws.Range("I60").Value = 53332.332
$ws.Range("J60").Value = 53436.332
$ws.Range("K60").Value = 53332.332
$ws.Range("L60").Value = 53436.332
$ws.Range("M60").Value = -52821.332
$ws.Range("N60").Value = -54458.332

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100.4
$ws.Range("J2").Value = 100.4
$ws.Range("L2").Value = 602.4000000000001
$ws.Range("N2").Value = -828.4000000000001
$ws.Range("H12").Value = 139.77777
$ws.Range("I12").Value = 133
$ws.Range("J12").Value = 143.16667
$ws.Range("K12").Value = 399
$ws.Range("L12").Value = 429.50001
$ws.Range("M12").Value = -226
$ws.Range("N12").Value = -775.50001
$ws.Range("H113").Value = 1356.5294
$ws.Range("I113").Value = 617.9048
$ws.Range("J113").Value = 2549.6924
$ws.Range("K113").Value = 1853.7144
$ws.Range("L113").Value = 7649.0772
$ws.Range("M113").Value = 316.2855999999999
$ws.Range("N113").Value = -11989.0772
$ws.Range("H131").Value = 1986.4
$ws.Range("J131").Value = 1983
$ws.Range("L131").Value = 5949
$ws.Range("N131").Value = -16029
$ws.Range("H139").Value = 4073.8
$ws.Range("I139").Value = 4073.8
$ws.Range("K139").Value = 12221.4
$ws.Range("M139").Value = -7081.400000000001
$ws.Range("H140").Value = 627174.9
$ws.Range("I140").Value = 627174.9
$ws.Range("K140").Value = 1881524.7
$ws.Range("M140").Value = -1876344.7
$ws.Range("H141").Value = 7560.75
$ws.Range("I141").Value = 7560.75
$ws.Range("K141").Value = 22682.25
$ws.Range("M141").Value = -17502.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9122.799999999999
$ws.Range("I70").Value = 8538
$ws.Range("K70").Value = 8538
$ws.Range("M70").Value = -8268
$ws.Range("H73").Value = 9122.799999999999
$ws.Range("I73").Value = 8538
$ws.Range("K73").Value = 8538
$ws.Range("M73").Value = -7602
$ws.Range("H122").Value = 4777.4
$ws.Range("I122").Value = 4513.5
$ws.Range("K122").Value = 13540.5
$ws.Range("M122").Value = -11090.5
$ws.Range("H126").Value = 4455.4443
$ws.Range("I126").Value = 4300
$ws.Range("K126").Value = 12900
$ws.Range("M126").Value = -10430
$ws.Range("H132").Value = 2166.1333
$ws.Range("I132").Value = 2184.4167
$ws.Range("K132").Value = 6553.250100000001
$ws.Range("M132").Value = -4023.250100000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3844.3572
$ws.Range("I7").Value = 3573.8
$ws.Range("J7").Value = 4520.75
$ws.Range("K7").Value = 3573.8
$ws.Range("L7").Value = 4520.75
$ws.Range("M7").Value = -3461.8
$ws.Range("N7").Value = -4744.75
$ws.Range("H24").Value = 25000
$ws.Range("J24").Value = 25000
$ws.Range("L24").Value = 25000
$ws.Range("N24").Value = -25686
$ws.Range("H40").Value = 4896.6
$ws.Range("I40").Value = 4896.6
$ws.Range("K40").Value = 4896.6
$ws.Range("M40").Value = -4760.6
$ws.Range("H82").Value = 1799.0834
$ws.Range("I82").Value = 1033
$ws.Range("J82").Value = 2054.4443
$ws.Range("K82").Value = 1033
$ws.Range("L82").Value = 2054.4443
$ws.Range("M82").Value = -672
$ws.Range("N82").Value = -2776.4443
$ws.Range("H85").Value = 1799.0834
$ws.Range("I85").Value = 1033
$ws.Range("J85").Value = 2054.4443
$ws.Range("K85").Value = 1033
$ws.Range("L85").Value = 2054.4443
$ws.Range("M85").Value = 215
$ws.Range("N85").Value = -4550.4443
$ws.Range("H100").Value = 1641.75
$ws.Range("I100").Value = 1608.6
$ws.Range("K100").Value = 1608.6
$ws.Range("M100").Value = -1067.6
$ws.Range("H126").Value = 3844.3572
$ws.Range("I126").Value = 3573.8
$ws.Range("J126").Value = 4520.75
$ws.Range("K126").Value = 10721.4
$ws.Range("L126").Value = 13562.25
$ws.Range("M126").Value = -8251.400000000001
$ws.Range("N126").Value = -18502.25
$ws.Range("H132").Value = 2045.1
$ws.Range("I132").Value = 2016.2273
$ws.Range("K132").Value = 6048.6819
$ws.Range("M132").Value = -3518.6819

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 30000
$ws.Range("I21").Value = 30000
$ws.Range("K21").Value = 30000
$ws.Range("M21").Value = -29765
$ws.Range("H35").Value = 30000
$ws.Range("I35").Value = 30000
$ws.Range("K35").Value = 30000
$ws.Range("M35").Value = -29710
$ws.Range("H81").Value = 2901.8096
$ws.Range("I81").Value = 2891.5264
$ws.Range("K81").Value = 5783.0528
$ws.Range("M81").Value = -4722.0528
$ws.Range("H84").Value = 2901.8096
$ws.Range("I84").Value = 2891.5264
$ws.Range("K84").Value = 28915.264
$ws.Range("M84").Value = -23611.264
$ws.Range("H126").Value = 5497.6665
$ws.Range("I126").Value = 2493
$ws.Range("K126").Value = 7479
$ws.Range("M126").Value = -5009
$ws.Range("H136").Value = 7179.2085
$ws.Range("I136").Value = 7179.2085
$ws.Range("K136").Value = 21537.6255
$ws.Range("M136").Value = -18987.6255
